$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.181.36"
$ws.Range("E2").Value = "  -3.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.815.71"
$ws.Range("E3").Value = "  +1.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.63"
$ws.Range("E5").Value = "  -4.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.74"
$ws.Range("E6").Value = "  -4.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.811.06"
$ws.Range("E7").Value = "  +1.27%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -1.13%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -6.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").Value = "  -0.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  -3.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.37"
$ws.Range("E13").Value = "  -5.61%  "

# Row 14
$ws.Range("E14").Value = "  -5.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.455.13"
$ws.Range("E15").Value = "  +1.49%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.821.59"
$ws.Range("E16").Value = "  +1.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.283.17"
$ws.Range("E17").Value = "  -3.06%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  -5.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.04"
$ws.Range("E20").Value = "  -3.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.57"
$ws.Range("E21").Value = "  -3.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.735"
$ws.Range("E23").Value = "  +1.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.93"
$ws.Range("E24").Value = "  -1.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  -7.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000137"
$ws.Range("E26").Value = "  +0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("E27").Value = "  -6.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  -9.81%  "

# Row 29
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("E30").Value = "  -0.16%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("E31").Value = "  -2.42%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.97"
$ws.Range("E32").Value = "  +7.20%  "

# Row 33
$ws.Range("E33").Value = "  -3.26%  "

# Row 34
$ws.Range("E34").Value = "  -4.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36
$ws.Range("E36").Value = "  -4.65%  "

# Row 37
$ws.Range("E37").Value = "  -2.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.80"
$ws.Range("E38").Value = "  -5.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.325"
$ws.Range("E39").Value = "  -7.31%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.77"
$ws.Range("E40").Value = "  +4.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.04"
$ws.Range("E41").Value = "  -1.96%  "

# Row 42
$ws.Range("E42").Value = "  -4.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  -10.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.30"
$ws.Range("E44").Value = "  -4.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.30"
$ws.Range("E45").Value = "  -7.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.849.94"
$ws.Range("E46").Value = "  -4.13%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0353"
$ws.Range("E48").Value = "  -3.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.03"
$ws.Range("E49").Value = "  +0.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.38"
$ws.Range("E50").Value = "  -3.75%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.35"
$ws.Range("E51").Value = "  +6.61%  "
